$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.334.27"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "2.272.10"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Formula = "'305.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Formula = "'97.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.76%  "

$ws.Range("D7").Formula = "'0.529"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Formula = "'0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").Formula = "'35.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.49%  "

$ws.Range("D11").Formula = "'0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("E12").Value = "  -2.50%  "

$ws.Range("D13").Formula = "'6.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "2.628.90"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").Formula = "'14.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").Value = "2.271.73"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Formula = "'0.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.06%  "

$ws.Range("D18").Value = "42.260.94"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Formula = "'12.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Formula = "'5.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").Formula = "'67.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").Formula = "'240.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").Formula = "'2.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").Formula = "'1.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Formula = "'23.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Formula = "'37.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.53%  "

$ws.Range("D29").Formula = "'9.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("D31").Formula = "'159.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Formula = "'5.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").Formula = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "

$ws.Range("D35").Formula = "'0.0742"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Formula = "'17.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").Formula = "'0.106"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Formula = "'1.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.53%  "

$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("D41").Formula = "'4.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.99%  "

$ws.Range("D42").Formula = "'2.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.47%  "

$ws.Range("D43").Value = "1.991.64"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Formula = "'0.0285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "'18.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("D46").Formula = "'2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.16%  "

$ws.Range("D47").Formula = "'9.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.97%  "

$ws.Range("D48").Formula = "'53.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("D49").Formula = "'1.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Formula = "'72.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").Formula = "'92.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
